$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F, shifting the old F (Labels) to G
$ws.Columns("F:F").Insert()

# Set the new column F header
$ws.Range("F1").Value = "Anzahl Linien"

# Populate the "Anzahl Linien" values for rows 2-100
$ws.Range("F2").Value = 3
$ws.Range("F3").Value = 5
$ws.Range("F4").Value = 5
$ws.Range("F5").Value = 8
$ws.Range("F6").Value = 8
$ws.Range("F7").Value = 16
$ws.Range("F8").Value = 18
$ws.Range("F9").Value = 7
$ws.Range("F10").Value = 7
$ws.Range("F11").Value = 11
$ws.Range("F12").Value = 11
$ws.Range("F13").Value = 8
$ws.Range("F14").Value = 3
$ws.Range("F15").Value = 9
$ws.Range("F16").Value = 4
$ws.Range("F17").Value = 8
$ws.Range("F18").Value = 4
$ws.Range("F19").Value = 6
$ws.Range("F20").Value = 5
$ws.Range("F21").Value = 7
$ws.Range("F22").Value = 14
$ws.Range("F23").Value = 12
$ws.Range("F24").Value = 10
$ws.Range("F25").Value = 11
$ws.Range("F26").Value = 5
$ws.Range("F27").Value = 5
$ws.Range("F28").Value = 5
$ws.Range("F29").Value = 9
$ws.Range("F30").Value = 4
$ws.Range("F31").Value = 8
$ws.Range("F32").Value = 6
$ws.Range("F33").Value = 4
$ws.Range("F34").Value = 6
$ws.Range("F35").Value = 15
$ws.Range("F36").Value = 14
$ws.Range("F37").Value = 6
$ws.Range("F38").Value = 6
$ws.Range("F39").Value = 3
$ws.Range("F40").Value = 3
$ws.Range("F41").Value = 5
$ws.Range("F42").Value = 8
$ws.Range("F43").Value = 15
$ws.Range("F44").Value = 6
$ws.Range("F45").Value = 9
$ws.Range("F46").Value = 8
$ws.Range("F47").Value = 11
$ws.Range("F48").Value = 10
$ws.Range("F49").Value = 6
$ws.Range("F50").Value = 10
$ws.Range("F51").Value = 16
$ws.Range("F52").Value = 13
$ws.Range("F53").Value = 12
$ws.Range("F54").Value = 17
$ws.Range("F55").Value = 9
$ws.Range("F56").Value = 12
$ws.Range("F57").Value = 4
$ws.Range("F58").Value = 6
$ws.Range("F59").Value = 3
$ws.Range("F60").Value = 5
$ws.Range("F61").Value = 7
$ws.Range("F62").Value = 13
$ws.Range("F63").Value = 15
$ws.Range("F64").Value = 13
$ws.Range("F65").Value = 5
$ws.Range("F66").Value = 10
$ws.Range("F67").Value = 11
$ws.Range("F68").Value = 13
$ws.Range("F69").Value = 6
$ws.Range("F70").Value = 5
$ws.Range("F71").Value = 6
$ws.Range("F72").Value = 9
$ws.Range("F73").Value = 10
$ws.Range("F74").Value = 3
$ws.Range("F75").Value = 7
$ws.Range("F76").Value = 9
$ws.Range("F77").Value = 5
$ws.Range("F78").Value = 10
$ws.Range("F79").Value = 4
$ws.Range("F80").Value = 2
$ws.Range("F81").Value = 2
$ws.Range("F82").Value = 4
$ws.Range("F83").Value = 5
$ws.Range("F84").Value = 3
$ws.Range("F85").Value = 1
$ws.Range("F86").Value = 4
$ws.Range("F87").Value = 1
$ws.Range("F88").Value = 4
$ws.Range("F89").Value = 1
$ws.Range("F90").Value = 5
$ws.Range("F91").Value = 2
$ws.Range("F92").Value = 1
$ws.Range("F93").Value = 5
$ws.Range("F94").Value = 5
$ws.Range("F95").Value = 3
$ws.Range("F96").Value = 5
$ws.Range("F97").Value = 3
$ws.Range("F98").Value = 3
$ws.Range("F99").Value = 5
$ws.Range("F100").Value = 3
